$d = $word.ActiveDocument

# Replace the original placeholder text with "1" (keeps the existing
# paragraph mark / hidden _GoBack bookmark attached right after it).
$d.Content.Find.Execute("rhtxsjfuyjtjdufyjdytfjhty", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

# Insert three more paragraphs ("2", "3", "4") right before the existing
# trailing paragraph mark (position 1, i.e. right after "1" but before the
# bookmark/paragraph end). Using Range.InsertAfter keeps the bookmark
# anchored to paragraph 1 instead of migrating it to the new paragraphs.
$r = $d.Range(1, 1)
$r.InsertAfter("`r2`r3`r4")

# Finally, append one more (empty) paragraph after "4" using Find/Replace
# with a "^p" marker instead of Range.InsertParagraphAfter/InsertAfter("`r"),
# since the latter would materialize a spurious empty run in the new
# paragraph. Find/Replace reuses the original run-less paragraph mark.
$d.Content.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "4^p", 2)
